$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.629.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.04%  "

$ws.Range("D3").Value = "'1.852.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.34%  "

$ws.Range("D4").Value = "'1.033"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +2.91%  "

$ws.Range("D5").Value = "'321.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.16%  "

$ws.Range("E6").Value = "  +2.76%  "

$ws.Range("D7").Value = "'0.4383"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.60%  "

$ws.Range("D8").Value = "'0.3757"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.58%  "

$ws.Range("D9").Value = "'0.07410"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.22%  "

$ws.Range("E10").Value = "  +1.01%  "

$ws.Range("E11").Value = "  +3.13%  "

$ws.Range("D12").Value = "'1.861.96"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.26%  "

$ws.Range("D13").Value = "'5.521"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.21%  "

$ws.Range("D14").Value = "'6.706"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.20%  "

$ws.Range("D15").Value = "'0.07191"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.17%  "

$ws.Range("D16").Value = "'83.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.21%  "

$ws.Range("E17").Value = "  +3.33%  "

$ws.Range("D18").Value = "'0.000009050"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.31%  "

$ws.Range("E19").Value = "  +2.87%  "

$ws.Range("D20").Value = "'15.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.34%  "

$ws.Range("D21").Value = "'27.642.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.88%  "

$ws.Range("D22").Value = "'5.271"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.33%  "

$ws.Range("D23").Value = "'11.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.54%  "

$ws.Range("D24").Value = "'2.067.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.29%  "

$ws.Range("D25").Value = "'157.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.70%  "

$ws.Range("E26").Value = "  +3.76%  "

$ws.Range("D27").Value = "'18.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.78%  "

$ws.Range("D28").Value = "'5.301"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.72%  "

$ws.Range("D29").Value = "'1.942"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.92%  "

$ws.Range("D30").Value = "'116.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.00%  "

$ws.Range("D31").Value = "'0.09082"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.54%  "

$ws.Range("D32").Value = "'1.208"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.33%  "

$ws.Range("D33").Value = "'0.7682"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.53%  "

$ws.Range("D34").Value = "'4.514"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.71%  "

$ws.Range("D35").Value = "'2.885"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.77%  "

$ws.Range("E36").Value = "  +2.66%  "

$ws.Range("D37").Value = "'1.156"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.59%  "

$ws.Range("D38").Value = "'0.01983"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.06%  "

$ws.Range("D39").Value = "'0.05284"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.28%  "

$ws.Range("D40").Value = "'0.5190"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.13%  "

$ws.Range("D41").Value = "'2.815"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.81%  "

$ws.Range("D42").Value = "'0.1675"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.63%  "

$ws.Range("D43").Value = "'6.749"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.16%  "

$ws.Range("D44").Value = "'8.597"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.82%  "

$ws.Range("D45").Value = "'109.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.43%  "

$ws.Range("D46").Value = "'10.62"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.45%  "

$ws.Range("E47").Value = "  +3.85%  "

$ws.Range("D48").Value = "'0.4665"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.09%  "

$ws.Range("D49").Value = "'0.06397"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.90%  "

$ws.Range("D50").Value = "'1.889"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.93%  "

$ws.Range("D51").Value = "'39.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.93%  "
